{"js": "// Apply the set of text replacements described by the diff:\n// the worksheet date and each division-problem cell.\nconst body = context.document.body;\n\n{\n  const results = body.search(\"2024-03-24 Sunday\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '2024-03-24 Sunday', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"2024-03-25 Monday\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"42\u00f77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '42\u00f77=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"82\u00f77=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"27\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '27\u00f79=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"29\u00f72=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"91\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '91\u00f72=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"65\u00f77=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"98\u00f76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '98\u00f76=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"48\u00f73=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"29\u00f77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '29\u00f77=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"38\u00f79=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"90\u00f78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '90\u00f78=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"37\u00f74=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"23\u00f75=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '23\u00f75=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"46\u00f76=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"89\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '89\u00f72=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"94\u00f72=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"25\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '25\u00f79=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"56\u00f77=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"55\u00f77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '55\u00f77=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"64\u00f73=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"41\u00f77=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '41\u00f77=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"10\u00f74=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"77\u00f78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '77\u00f78=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"53\u00f72=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"40\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '40\u00f72=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"85\u00f72=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"50\u00f78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '50\u00f78=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"99\u00f78=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"41\u00f74=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '41\u00f74=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"22\u00f78=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"20\u00f79=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '20\u00f79=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"39\u00f73=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"32\u00f75=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '32\u00f75=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"34\u00f75=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"56\u00f78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '56\u00f78=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"22\u00f78=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"34\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '34\u00f72=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"60\u00f76=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"46\u00f73=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '46\u00f73=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"75\u00f74=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"33\u00f76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '33\u00f76=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"63\u00f79=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"77\u00f76=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '77\u00f76=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"69\u00f78=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"66\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '66\u00f72=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"18\u00f79=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"28\u00f78=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '28\u00f78=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"53\u00f77=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n{\n  const results = body.search(\"31\u00f72=\", { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '31\u00f72=', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"63\u00f76=\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\nreturn \"ok\";", "ps1": "# Apply the text replacements described by the diff: the worksheet\n# date and each two-digit division problem, one Find/Replace per cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-03-24 Sunday\", \"2024-03-25 Monday\"),\n  @(\"42\u00f77=\", \"82\u00f77=\"),\n  @(\"27\u00f79=\", \"29\u00f72=\"),\n  @(\"91\u00f72=\", \"65\u00f77=\"),\n  @(\"98\u00f76=\", \"48\u00f73=\"),\n  @(\"29\u00f77=\", \"38\u00f79=\"),\n  @(\"90\u00f78=\", \"37\u00f74=\"),\n  @(\"23\u00f75=\", \"46\u00f76=\"),\n  @(\"89\u00f72=\", \"94\u00f72=\"),\n  @(\"25\u00f79=\", \"56\u00f77=\"),\n  @(\"55\u00f77=\", \"64\u00f73=\"),\n  @(\"41\u00f77=\", \"10\u00f74=\"),\n  @(\"77\u00f78=\", \"53\u00f72=\"),\n  @(\"40\u00f72=\", \"85\u00f72=\"),\n  @(\"50\u00f78=\", \"99\u00f78=\"),\n  @(\"41\u00f74=\", \"22\u00f78=\"),\n  @(\"20\u00f79=\", \"39\u00f73=\"),\n  @(\"32\u00f75=\", \"34\u00f75=\"),\n  @(\"56\u00f78=\", \"22\u00f78=\"),\n  @(\"34\u00f72=\", \"60\u00f76=\"),\n  @(\"46\u00f73=\", \"75\u00f74=\"),\n  @(\"33\u00f76=\", \"63\u00f79=\"),\n  @(\"77\u00f76=\", \"69\u00f78=\"),\n  @(\"66\u00f72=\", \"18\u00f79=\"),\n  @(\"28\u00f78=\", \"53\u00f77=\"),\n  @(\"31\u00f72=\", \"63\u00f76=\"),\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $r = $d.Content\n  $r.Find.ClearFormatting()\n  $r.Find.Replacement.ClearFormatting()\n  $found = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n  if (-not $found) {\n    throw \"Find/Replace did not match: $oldText\"\n  }\n}\n\nWrite-Output \"done\""}
